# Natmi following Dr Hou advice
# Update the LR-pair summary sheet so the Sending cluster x Target cluster
# combinations form the complete 3x3 matrix (ECs, FAPs, sCs) instead of the
# previous matrix that skipped the "same cluster" (self) pairs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array corresponds to one data row (columns A..T):
# A  Sending cluster
# B  Ligand symbol
# C  Receptor symbol
# D  Target cluster
# E  Ligand-expressing cells
# F  Ligand detection rate
# G  Ligand average expression value
# H  Ligand total expression value
# I  Ligand derived specificity of average expression value
# J  Ligand derived specificity of total expression value
# K  Receptor-expressing cells
# L  Receptor detection rate
# M  Receptor average expression value
# N  Receptor total expression value
# O  Receptor derived specificity of average expression value
# P  Receptor derived specificity of total expression value
# Q  Edge average expression weight
# R  Edge total expression weight
# S  Edge average expression derived specificity
# T  Edge total expression derived specificity

$rows = @(
    @("ECs",  "Efnb2", "Epha4", "ECs",  3, 1, 36.899643,   110.698929,   0.7238945645409351, 0.7238945645409351, 2, 0.6666666666666666, 6.708176333333333, 20.124529,          0.4356329228871633, 0.4356329228871633, 247.529311881049,  2227.763806929441, 0.3153523050130979,  0.3153523050130979),
    @("ECs",  "Efnb2", "Epha4", "FAPs", 3, 1, 36.899643,   110.698929,   0.7238945645409351, 0.7238945645409351, 3, 1,                   6.789877333333333, 20.369632,          0.4409386339573907, 0.4409386339573907, 250.544049613792,  2254.896446524128, 0.3191930804178602,  0.3191930804178602),
    @("ECs",  "Efnb2", "Epha4", "sCs",  3, 1, 36.899643,   110.698929,   0.7238945645409351, 0.7238945645409351, 3, 1,                   1.900636333333334, 5.701909000000001, 0.1234284431554459, 0.1234284431554459, 70.132802172829,   631.195219555461, 0.0893491791099771,  0.0893491791099771),
    @("FAPs", "Efnb2", "Epha4", "ECs",  3, 1, 3.374819,    10.124457,    0.0662069584361419, 0.0662069584361419, 2, 0.6666666666666666, 6.708176333333333, 20.124529,          0.4356329228871633, 0.4356329228871633, 22.63888094508367, 203.749928505753, 0.02884193081900543, 0.02884193081900543),
    @("FAPs", "Efnb2", "Epha4", "FAPs", 3, 1, 3.374819,    10.124457,    0.0662069584361419, 0.0662069584361419, 3, 1,                   6.789877333333333, 20.369632,          0.4409386339573907, 0.4409386339573907, 22.91460703220267, 206.231463289824, 0.02919320581130615, 0.02919320581130615),
    @("FAPs", "Efnb2", "Epha4", "sCs",  3, 1, 3.374819,    10.124457,    0.0662069584361419, 0.0662069584361419, 3, 1,                   1.900636333333334, 5.701909000000001, 0.1234284431554459, 0.1234284431554459, 6.414303609823667, 57.728732488413,  0.008171821805830311, 0.008171821805830311),
    @("sCs",  "Efnb2", "Epha4", "ECs",  3, 1, 10.699319,   32.097957,    0.2098984770229228, 0.2098984770229228, 2, 0.6666666666666666, 6.708176333333333, 20.124529,          0.4356329228871633, 0.4356329228871633, 71.77291849858368, 645.9562664872529, 0.09143868705505996, 0.09143868705505996),
    @("sCs",  "Efnb2", "Epha4", "FAPs", 3, 1, 10.699319,   32.097957,    0.2098984770229228, 0.2098984770229228, 3, 1,                   6.789877333333333, 20.369632,          0.4409386339573907, 0.4409386339573907, 72.64706356020267, 653.823572041824,  0.09255234772822436, 0.09255234772822436),
    @("sCs",  "Efnb2", "Epha4", "sCs",  3, 1, 10.699319,   32.097957,    0.2098984770229228, 0.2098984770229228, 3, 1,                   1.900636333333334, 5.701909000000001, 0.1234284431554459, 0.1234284431554459, 20.33551443332367, 183.019629899913,  0.0259074422396385,  0.0259074422396385)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
